# Auto commit at 2025-10-22  7:59:08.96
# Append two new daily summary rows (102, 103) for date 45951 (2025-10-21),
# one for each station ("四方坪站" / "高岭站"), mirroring the formula
# pattern already used for the preceding days' rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---- Row 102: 四方坪站 ------------------------------------------------
$ws.Range("A102").Value = 45951
$ws.Range("B102").Value = "四方坪站"
$ws.Range("C102").Formula = "=18469/126"
$ws.Range("D102").Formula = "=C102/(24*60)"
$ws.Range("E102").Formula = "=9719.94/126"
$ws.Range("F102").Formula = "=3449.26/126"
$ws.Range("G102").Formula = "=9719.94/(18469/60)"
$ws.Range("H102").Formula = "=408/126"

# ---- Row 103: 高岭站 --------------------------------------------------
$ws.Range("A103").Value = 45951
$ws.Range("B103").Value = "高岭站"
$ws.Range("C103").Formula = "=9025/36"
$ws.Range("D103").Formula = "=C103/(24*60)"
$ws.Range("E103").Formula = "=5947.02/36"
$ws.Range("F103").Formula = "=1572.22/36"
$ws.Range("G103").Formula = "=5947.02/(9025/60)"
$ws.Range("H103").Formula = "=209/36"

# ---- Scroll / selection bookkeeping (mirrors author's view state) ----
$excel.ActiveWindow.ScrollRow = 88
$ws.Range("I98").Select() | Out-Null
